# Cover letter revision for new submission.
# - Normalizes paragraph styles to "Normal" (explicit style application).
# - Merges the split "Drug-Drug" run sequence into a single run.
# - Drops the stale lastRenderedPageBreak marker ahead of "We trust...".
# - Tightens up page setup (header/footer distance, page-number format).

$d = $word.ActiveDocument

# --- 1. Merge the three runs of the manuscript-title paragraph into one run ---
$mergedTitle = 'I am submitting the manuscript titled "SNF-CNN: Predicting Comprehensive Drug-Drug Interaction via Similarity Network Fusion and Convolutional Neural Networks" for consideration for publication in the Journal of Bioinformatics.'

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("I am submitting the manuscript titled")) {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "TEMP_PLACEHOLDER_MERGE"
        $p2 = $p
        $r2 = $p2.Range
        $r2.End = $r2.End - 1
        $r2.Text = $mergedTitle
        break
    }
}

# --- 2. Strip the leftover lastRenderedPageBreak on the "We trust..." paragraph ---
$trustText = "We trust that our research will be of interest to the readership of the Journal of Bioinformatics and contribute to the ongoing discourse in the field of drug-drug interactions. We look forward to the opportunity for our work to be peer-reviewed and, hopefully, published in your esteemed journal."

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("We trust that our research")) {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "TEMP_PLACEHOLDER_BREAK"
        $p2 = $p
        $r2 = $p2.Range
        $r2.End = $r2.End - 1
        $r2.Text = $trustText
        break
    }
}

# --- 3. Explicitly (re)apply the "Normal" paragraph style to every paragraph ---
foreach ($p in $d.Paragraphs) {
    $p.Style = "Normal"
}

# --- 3b. Re-affirm the resolved paragraph formatting on the closing signature line ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.WidowControl = $true
$lastPara.ReadingOrder = 0
$lastPara.LineSpacingRule = 5
$lastPara.LineSpacing = 12.95
$lastPara.SpaceBefore = 0
$lastPara.SpaceAfter = 8

# --- 4. Page setup: headers/footers flush (0 distance) and decimal page numbers ---
$section = $d.Sections(1)
$ps = $section.PageSetup
$ps.HeaderDistance = 0
$ps.FooterDistance = 0

$hf = $section.Headers(1)
$hf.PageNumbers.NumberStyle = 0

$section.ProtectedForForms = $false

Write-Output "done"
